# Update cryptocurrency price/volume figures to the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.284.70"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.664.62"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.83%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5331"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2649"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.564"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "1.662.42"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "1.892.35"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "0.0₅8214"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.685"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.036"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.199"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.481"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05872"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.280"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.609"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9629"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.824"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.416"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5807"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01609"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8652"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.859"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").Value = "1.050.64"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.009"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").Value = "1.803.08"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  -5.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4383"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.058"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05165"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.421"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.76%  "
